# Insert a new column before column C ("AV Delay"), shifting the
# existing Atrial Amplitude..ARP columns one to the right, then fix up
# the handful of row-2 values that differ from a pure shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:H -> D:I and open up column C for the new field.
$ws.Columns("C:C").Insert()

# New header + value for the inserted "AV Delay" column.
$ws.Range("C1").Value = "AV Delay"
$ws.Range("C2").Value = "'70"

# These values are stored as text in the sheet (same as every other
# cell here), so force-text them with a leading apostrophe to avoid
# Excel reinterpreting the numeric-looking strings as numbers.
$ws.Range("E2").Value = "'1.2"
$ws.Range("G2").Value = "'1.2"
$ws.Range("H2").Value = "'150"
$ws.Range("I2").Value = "'150"
